$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 4 (old Resolving-Mac -> Resolving-Mac self loop row) first,
# so remaining data only spans rows 1-3.
$ws.Rows.Item(4).Delete()

# Row 2: MuSCs -> Sectm1a/Cd7 -> FAPs (updated TPM-derived values)
$ws.Cells.Item(2, 1).Value2 = "MuSCs"
$ws.Cells.Item(2, 2).Value2 = "Sectm1a"
$ws.Cells.Item(2, 3).Value2 = "Cd7"
$ws.Cells.Item(2, 4).Value2 = "FAPs"
$ws.Cells.Item(2, 5).Value2 = 1
$ws.Cells.Item(2, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 7).Value2 = 0.124553
$ws.Cells.Item(2, 8).Value2 = 0.373659
$ws.Cells.Item(2, 9).Value2 = 1
$ws.Cells.Item(2, 10).Value2 = 1
$ws.Cells.Item(2, 11).Value2 = 1
$ws.Cells.Item(2, 12).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 13).Value2 = 0.5011856666666666
$ws.Cells.Item(2, 14).Value2 = 1.503557
$ws.Cells.Item(2, 15).Value2 = 0.5449482835232878
$ws.Cells.Item(2, 16).Value2 = 0.5449482835232877
$ws.Cells.Item(2, 17).Value2 = 0.06242417834033333
$ws.Cells.Item(2, 18).Value2 = 0.561817605063
$ws.Cells.Item(2, 19).Value2 = 0.5449482835232878
$ws.Cells.Item(2, 20).Value2 = 0.5449482835232877

# Row 3: MuSCs -> Sectm1a/Cd7 -> Resolving-Mac (updated TPM-derived values)
$ws.Cells.Item(3, 1).Value2 = "MuSCs"
$ws.Cells.Item(3, 2).Value2 = "Sectm1a"
$ws.Cells.Item(3, 3).Value2 = "Cd7"
$ws.Cells.Item(3, 4).Value2 = "Resolving-Mac"
$ws.Cells.Item(3, 5).Value2 = 1
$ws.Cells.Item(3, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(3, 7).Value2 = 0.124553
$ws.Cells.Item(3, 8).Value2 = 0.373659
$ws.Cells.Item(3, 9).Value2 = 1
$ws.Cells.Item(3, 10).Value2 = 1
$ws.Cells.Item(3, 11).Value2 = 2
$ws.Cells.Item(3, 12).Value2 = 0.6666666666666666
$ws.Cells.Item(3, 13).Value2 = 0.4185083333333333
$ws.Cells.Item(3, 14).Value2 = 1.255525
$ws.Cells.Item(3, 15).Value2 = 0.4550517164767122
$ws.Cells.Item(3, 16).Value2 = 0.4550517164767122
$ws.Cells.Item(3, 17).Value2 = 0.05212646844166667
$ws.Cells.Item(3, 18).Value2 = 0.4691382159750001
$ws.Cells.Item(3, 19).Value2 = 0.4550517164767122
$ws.Cells.Item(3, 20).Value2 = 0.4550517164767122
